$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ticket #121 as a new row 14. Columns A ("121") and G
# ("2025-05-06") look like a number / a date respectively, so a plain
# assignment would get auto-typed by Excel. The rest of the sheet
# stores every value as literal text, so prefix those two with an
# apostrophe (the standard "force text" trick) to keep them text too.
$ws.Range("A14").Value = "'121"
$ws.Range("B14").Value = "SAP Hana"
$ws.Range("C14").Value = "ramya"
$ws.Range("D14").Value = "Login being denied"
$ws.Range("E14").Value = "Trying to login but my credentials are denied."
$ws.Range("F14").Value = "Thats becoz some caps issue in login credentials . pls try andain and let me know if it persists."
$ws.Range("G14").Value = "'2025-05-06"
